$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53
$ws.Cells.Item(53, 1).Value = 112230608
$ws.Cells.Item(53, 2).Value = 99413
$ws.Cells.Item(53, 3).Value = "Ovaliderad"
$ws.Cells.Item(53, 4).Value = "LC"
$ws.Cells.Item(53, 5).Value = 221235
$ws.Cells.Item(53, 6).Value = "Vårärt"
$ws.Cells.Item(53, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(53, 8).Value = "(L.) Bernh."
$ws.Cells.Item(53, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(53, 17).Value = 571931
$ws.Cells.Item(53, 18).Value = 6697694
$ws.Cells.Item(53, 19).Value = 15
$ws.Cells.Item(53, 20).Value = "Dalarna"
$ws.Cells.Item(53, 21).Value = "Hedemora"
$ws.Cells.Item(53, 22).Value = "Dalarna"
$ws.Cells.Item(53, 23).Value = "Husby"
$ws.Cells.Item(53, 25).Value = "'2023-09-21"
$ws.Cells.Item(53, 27).Value = "'2023-09-21"
$ws.Cells.Item(53, 30).Value = $false
$ws.Cells.Item(53, 31).Value = $false
$ws.Cells.Item(53, 33).Value = $false
$ws.Cells.Item(53, 49).Value = "Philipp Weiss"
$ws.Cells.Item(53, 50).Value = "Philipp Weiss"

# Row 54
$ws.Cells.Item(54, 1).Value = 112230611
$ws.Cells.Item(54, 2).Value = 4711
$ws.Cells.Item(54, 3).Value = "Ovaliderad"
$ws.Cells.Item(54, 4).Value = "LC"
$ws.Cells.Item(54, 5).Value = 100299
$ws.Cells.Item(54, 6).Value = "Thomsons trägnagare"
$ws.Cells.Item(54, 7).Value = "Cacotemnus thomsoni"
$ws.Cells.Item(54, 8).Value = "(Kraatz, 1881)"
$ws.Cells.Item(54, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(54, 17).Value = 571834
$ws.Cells.Item(54, 18).Value = 6697641
$ws.Cells.Item(54, 19).Value = 15
$ws.Cells.Item(54, 20).Value = "Dalarna"
$ws.Cells.Item(54, 21).Value = "Hedemora"
$ws.Cells.Item(54, 22).Value = "Dalarna"
$ws.Cells.Item(54, 23).Value = "Husby"
$ws.Cells.Item(54, 25).Value = "'2023-09-21"
$ws.Cells.Item(54, 27).Value = "'2023-09-21"
$ws.Cells.Item(54, 30).Value = $false
$ws.Cells.Item(54, 31).Value = $false
$ws.Cells.Item(54, 33).Value = $false
$ws.Cells.Item(54, 49).Value = "Philipp Weiss"
$ws.Cells.Item(54, 50).Value = "Philipp Weiss"

# Row 55
$ws.Cells.Item(55, 1).Value = 112230613
$ws.Cells.Item(55, 2).Value = 89405
$ws.Cells.Item(55, 3).Value = "Ovaliderad"
$ws.Cells.Item(55, 4).Value = "NT"
$ws.Cells.Item(55, 5).Value = 1202
$ws.Cells.Item(55, 6).Value = "Ullticka"
$ws.Cells.Item(55, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(55, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(55, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(55, 17).Value = 571799
$ws.Cells.Item(55, 18).Value = 6697620
$ws.Cells.Item(55, 19).Value = 15
$ws.Cells.Item(55, 20).Value = "Dalarna"
$ws.Cells.Item(55, 21).Value = "Hedemora"
$ws.Cells.Item(55, 22).Value = "Dalarna"
$ws.Cells.Item(55, 23).Value = "Husby"
$ws.Cells.Item(55, 25).Value = "'2023-09-21"
$ws.Cells.Item(55, 27).Value = "'2023-09-21"
$ws.Cells.Item(55, 30).Value = $false
$ws.Cells.Item(55, 31).Value = $false
$ws.Cells.Item(55, 33).Value = $false
$ws.Cells.Item(55, 49).Value = "Philipp Weiss"
$ws.Cells.Item(55, 50).Value = "Philipp Weiss"

# Row 56
$ws.Cells.Item(56, 1).Value = 112230606
$ws.Cells.Item(56, 2).Value = 56543
$ws.Cells.Item(56, 3).Value = "Ovaliderad"
$ws.Cells.Item(56, 4).Value = "NT"
$ws.Cells.Item(56, 5).Value = 103021
$ws.Cells.Item(56, 6).Value = "Talltita"
$ws.Cells.Item(56, 7).Value = "Poecile montanus"
$ws.Cells.Item(56, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(56, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(56, 17).Value = 571961
$ws.Cells.Item(56, 18).Value = 6697705
$ws.Cells.Item(56, 19).Value = 15
$ws.Cells.Item(56, 20).Value = "Dalarna"
$ws.Cells.Item(56, 21).Value = "Hedemora"
$ws.Cells.Item(56, 22).Value = "Dalarna"
$ws.Cells.Item(56, 23).Value = "Husby"
$ws.Cells.Item(56, 25).Value = "'2023-09-21"
$ws.Cells.Item(56, 27).Value = "'2023-09-21"
$ws.Cells.Item(56, 30).Value = $false
$ws.Cells.Item(56, 31).Value = $false
$ws.Cells.Item(56, 33).Value = $false
$ws.Cells.Item(56, 49).Value = "Philipp Weiss"
$ws.Cells.Item(56, 50).Value = "Philipp Weiss"

# Row 57
$ws.Cells.Item(57, 1).Value = 112230614
$ws.Cells.Item(57, 2).Value = 78512
$ws.Cells.Item(57, 3).Value = "Ovaliderad"
$ws.Cells.Item(57, 4).Value = "LC"
$ws.Cells.Item(57, 5).Value = 6456
$ws.Cells.Item(57, 6).Value = "Skinnlav"
$ws.Cells.Item(57, 7).Value = "Leptogium saturninum"
$ws.Cells.Item(57, 8).Value = "(Dicks.) Nyl."
$ws.Cells.Item(57, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(57, 17).Value = 571792
$ws.Cells.Item(57, 18).Value = 6697651
$ws.Cells.Item(57, 19).Value = 15
$ws.Cells.Item(57, 20).Value = "Dalarna"
$ws.Cells.Item(57, 21).Value = "Hedemora"
$ws.Cells.Item(57, 22).Value = "Dalarna"
$ws.Cells.Item(57, 23).Value = "Husby"
$ws.Cells.Item(57, 25).Value = "'2023-09-21"
$ws.Cells.Item(57, 27).Value = "'2023-09-21"
$ws.Cells.Item(57, 30).Value = $false
$ws.Cells.Item(57, 31).Value = $false
$ws.Cells.Item(57, 33).Value = $false
$ws.Cells.Item(57, 49).Value = "Philipp Weiss"
$ws.Cells.Item(57, 50).Value = "Philipp Weiss"

# Row 58
$ws.Cells.Item(58, 1).Value = 112230605
$ws.Cells.Item(58, 2).Value = 99413
$ws.Cells.Item(58, 3).Value = "Ovaliderad"
$ws.Cells.Item(58, 4).Value = "LC"
$ws.Cells.Item(58, 5).Value = 221235
$ws.Cells.Item(58, 6).Value = "Vårärt"
$ws.Cells.Item(58, 7).Value = "Lathyrus vernus"
$ws.Cells.Item(58, 8).Value = "(L.) Bernh."
$ws.Cells.Item(58, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(58, 17).Value = 571995
$ws.Cells.Item(58, 18).Value = 6697876
$ws.Cells.Item(58, 19).Value = 15
$ws.Cells.Item(58, 20).Value = "Dalarna"
$ws.Cells.Item(58, 21).Value = "Hedemora"
$ws.Cells.Item(58, 22).Value = "Dalarna"
$ws.Cells.Item(58, 23).Value = "Husby"
$ws.Cells.Item(58, 25).Value = "'2023-09-21"
$ws.Cells.Item(58, 27).Value = "'2023-09-21"
$ws.Cells.Item(58, 30).Value = $false
$ws.Cells.Item(58, 31).Value = $false
$ws.Cells.Item(58, 33).Value = $false
$ws.Cells.Item(58, 49).Value = "Philipp Weiss"
$ws.Cells.Item(58, 50).Value = "Philipp Weiss"

# Row 59
$ws.Cells.Item(59, 1).Value = 112230603
$ws.Cells.Item(59, 2).Value = 78512
$ws.Cells.Item(59, 3).Value = "Ovaliderad"
$ws.Cells.Item(59, 4).Value = "LC"
$ws.Cells.Item(59, 5).Value = 6456
$ws.Cells.Item(59, 6).Value = "Skinnlav"
$ws.Cells.Item(59, 7).Value = "Leptogium saturninum"
$ws.Cells.Item(59, 8).Value = "(Dicks.) Nyl."
$ws.Cells.Item(59, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(59, 17).Value = 572018
$ws.Cells.Item(59, 18).Value = 6697738
$ws.Cells.Item(59, 19).Value = 15
$ws.Cells.Item(59, 20).Value = "Dalarna"
$ws.Cells.Item(59, 21).Value = "Hedemora"
$ws.Cells.Item(59, 22).Value = "Dalarna"
$ws.Cells.Item(59, 23).Value = "Husby"
$ws.Cells.Item(59, 25).Value = "'2023-09-21"
$ws.Cells.Item(59, 27).Value = "'2023-09-21"
$ws.Cells.Item(59, 30).Value = $false
$ws.Cells.Item(59, 31).Value = $false
$ws.Cells.Item(59, 33).Value = $false
$ws.Cells.Item(59, 49).Value = "Philipp Weiss"
$ws.Cells.Item(59, 50).Value = "Philipp Weiss"

# Row 60
$ws.Cells.Item(60, 1).Value = 112230612
$ws.Cells.Item(60, 2).Value = 12274
$ws.Cells.Item(60, 3).Value = "Ovaliderad"
$ws.Cells.Item(60, 4).Value = "NT"
$ws.Cells.Item(60, 5).Value = 102016
$ws.Cells.Item(60, 6).Value = "Gropig brunbagge"
$ws.Cells.Item(60, 7).Value = "Zilora ferruginea"
$ws.Cells.Item(60, 8).Value = "(Paykull, 1798)"
$ws.Cells.Item(60, 11).Value = "larv/nymf"
$ws.Cells.Item(60, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(60, 17).Value = 571800
$ws.Cells.Item(60, 18).Value = 6697623
$ws.Cells.Item(60, 19).Value = 15
$ws.Cells.Item(60, 20).Value = "Dalarna"
$ws.Cells.Item(60, 21).Value = "Hedemora"
$ws.Cells.Item(60, 22).Value = "Dalarna"
$ws.Cells.Item(60, 23).Value = "Husby"
$ws.Cells.Item(60, 25).Value = "'2023-09-21"
$ws.Cells.Item(60, 27).Value = "'2023-09-21"
$ws.Cells.Item(60, 30).Value = $false
$ws.Cells.Item(60, 31).Value = $false
$ws.Cells.Item(60, 33).Value = $false
$ws.Cells.Item(60, 49).Value = "Philipp Weiss"
$ws.Cells.Item(60, 50).Value = "Philipp Weiss"

# Row 61
$ws.Cells.Item(61, 1).Value = 112230610
$ws.Cells.Item(61, 2).Value = 90332
$ws.Cells.Item(61, 3).Value = "Ovaliderad"
$ws.Cells.Item(61, 4).Value = "LC"
$ws.Cells.Item(61, 5).Value = 4769
$ws.Cells.Item(61, 6).Value = "Svavelriska"
$ws.Cells.Item(61, 7).Value = "Lactarius scrobiculatus"
$ws.Cells.Item(61, 8).Value = "(Scop.:Fr.) Fr."
$ws.Cells.Item(61, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(61, 17).Value = 571853
$ws.Cells.Item(61, 18).Value = 6697760
$ws.Cells.Item(61, 19).Value = 15
$ws.Cells.Item(61, 20).Value = "Dalarna"
$ws.Cells.Item(61, 21).Value = "Hedemora"
$ws.Cells.Item(61, 22).Value = "Dalarna"
$ws.Cells.Item(61, 23).Value = "Husby"
$ws.Cells.Item(61, 25).Value = "'2023-09-21"
$ws.Cells.Item(61, 27).Value = "'2023-09-21"
$ws.Cells.Item(61, 30).Value = $false
$ws.Cells.Item(61, 31).Value = $false
$ws.Cells.Item(61, 33).Value = $false
$ws.Cells.Item(61, 49).Value = "Philipp Weiss"
$ws.Cells.Item(61, 50).Value = "Philipp Weiss"

# Row 62
$ws.Cells.Item(62, 1).Value = 112230604
$ws.Cells.Item(62, 2).Value = 101703
$ws.Cells.Item(62, 3).Value = "Ovaliderad"
$ws.Cells.Item(62, 4).Value = "LC"
$ws.Cells.Item(62, 5).Value = 222412
$ws.Cells.Item(62, 6).Value = "Tibast"
$ws.Cells.Item(62, 7).Value = "Daphne mezereum"
$ws.Cells.Item(62, 8).Value = "L."
$ws.Cells.Item(62, 16).Value = "Nordbäcksbo, Dlr"
$ws.Cells.Item(62, 17).Value = 571996
$ws.Cells.Item(62, 18).Value = 6697876
$ws.Cells.Item(62, 19).Value = 15
$ws.Cells.Item(62, 20).Value = "Dalarna"
$ws.Cells.Item(62, 21).Value = "Hedemora"
$ws.Cells.Item(62, 22).Value = "Dalarna"
$ws.Cells.Item(62, 23).Value = "Husby"
$ws.Cells.Item(62, 25).Value = "'2023-09-21"
$ws.Cells.Item(62, 27).Value = "'2023-09-21"
$ws.Cells.Item(62, 30).Value = $false
$ws.Cells.Item(62, 31).Value = $false
$ws.Cells.Item(62, 33).Value = $false
$ws.Cells.Item(62, 49).Value = "Philipp Weiss"
$ws.Cells.Item(62, 50).Value = "Philipp Weiss"
